$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data (daily time-series ending 2026/01/05 .. 2027/01/05) gained one
# more reading: a new "2026/01/05 月 23:00" row that belongs right before the
# existing "2026/12/29" block, at row 573. Inserting a whole row there shifts
# every following row down by one (573-614 -> 574-615), which is exactly what
# the diff shows - the old last row (2027/01/05 火 7) simply becomes row 615
# with no edits needed.
$ws.Rows.Item(573).Insert()

# Copy the date/weekday text from the row just above (same day: 2026/01/05,
# 月) instead of typing the literal "2026/01/05" into .Value, since Excel
# would otherwise auto-convert a date-shaped string into a date serial.
$ws.Range("A572:B572").Copy($ws.Range("A573:B573"))

# Fill in the numeric columns for the newly inserted row.
$ws.Range("C573").Value = 23
$ws.Range("D573").Value = 201
